$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.892.40"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "2.612.23"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'582.09"
$ws.Range("E5").Value = "  +4.69%  "
$ws.Range("D6").Value = "'144.28"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "2.638.57"
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("D13").Value = "'0.373"
$ws.Range("E13").Value = "  +6.45%  "
$ws.Range("D14").Value = "3.080.05"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").Value = "60.857.93"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "'23.44"
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("E17").Value = "  +4.44%  "
$ws.Range("D18").Value = "2.626.12"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "'11.31"
$ws.Range("E19").Value = "  +9.64%  "
$ws.Range("D20").Value = "'4.68"
$ws.Range("E20").Value = "  +3.22%  "
$ws.Range("D21").Value = "'350.36"
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("D22").Value = "'6.97"
$ws.Range("E22").Value = "  +7.93%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "'0.517"
$ws.Range("E24").Value = "  +8.12%  "
$ws.Range("D25").Value = "'63.27"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D28").Value = "'7.94"
$ws.Range("E28").Value = "  +7.73%  "
$ws.Range("D29").Value = "0.0₃0801"
$ws.Range("E29").Value = "  +3.77%  "
$ws.Range("E30").Value = "  +9.16%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.35"
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("D33").Value = "'162.92"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").Value = "'19.60"
$ws.Range("E34").Value = "  +2.89%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.32"
$ws.Range("E35").Value = "  +5.95%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = "  +13.69%  "
$ws.Range("E37").Value = "  +6.55%  "
$ws.Range("E38").Value = "  +10.36%  "
$ws.Range("D39").Value = "'37.98"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("D40").Value = "'3.92"
$ws.Range("E40").Value = "  +6.97%  "
$ws.Range("D41").Value = "'309.49"
$ws.Range("E41").Value = "  +7.05%  "
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "'134.36"
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("D44").Value = "'20.42"
$ws.Range("E44").Value = "  +9.58%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.610"
$ws.Range("E45").Value = "  +3.00%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'0.995"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'5.04"
$ws.Range("E47").Value = "  +11.75%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.0984"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'19.88"
$ws.Range("E49").Value = "  +5.41%  "
$ws.Range("E50").Value = "  +4.21%  "
$ws.Range("E51").Value = "  +4.08%  "
